$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.404.07'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').Value = '2.279.70'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  -0.74%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.94'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.48'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.28'
$ws.Range('E10').Value = '  -2.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0796'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.26'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '2.621.43'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '2.276.21'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.814'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.68'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '45.219.82'
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.25'
$ws.Range('E19').Value = '  +13.32%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.05'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.51'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.22'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.89'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '41.35'
$ws.Range('E27').Value = '  +11.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.71'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '151.56'
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.55'
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0793'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.58'
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.94'
$ws.Range('E35').Value = '  -6.48%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.118'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value = '  -6.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.96'
$ws.Range('E39').Value = '  +5.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0313'
$ws.Range('E40').Value = '  +6.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.25'
$ws.Range('E41').Value = '  -3.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.56'
$ws.Range('E42').Value = '  -9.44%  '
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.96'
$ws.Range('E44').Value = '  +11.55%  '
$ws.Range('D45').Value = '1.768.53'
$ws.Range('E45').Value = '  -2.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.193'
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '70.35'
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.52'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.74'
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.85'
$ws.Range('E51').Value = '  +0.19%  '
